$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 69: item_recipe_hair_of_samson (recipe row)
$ws.Cells.Item(69, 1).Value = "item_recipe_hair_of_samson"
$ws.Cells.Item(69, 2).Value = -1
$ws.Cells.Item(69, 3).Value = 3074

# Row 70: item_hair_of_samson (actual item row)
$ws.Cells.Item(70, 1).Value = "item_hair_of_samson"
$ws.Cells.Item(70, 2).Value = 30
$ws.Cells.Item(70, 3).Value = 3075
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 6
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0
$ws.Cells.Item(70, 14).Value = 0
$ws.Cells.Item(70, 15).Value = 0
$ws.Cells.Item(70, 16).Value = 0
$ws.Cells.Item(70, 17).Value = 0
$ws.Cells.Item(70, 18).Value = 0
$ws.Cells.Item(70, 19).Value = 0
$ws.Cells.Item(70, 20).Value = 0
$ws.Cells.Item(70, 21).Value = 0
$ws.Cells.Item(70, 22).Value = 0
$ws.Cells.Item(70, 23).Value = 0
$ws.Cells.Item(70, 24).Value = 0
$ws.Cells.Item(70, 25).Value = 0
$ws.Cells.Item(70, 26).Value = 0

# AC70 - shop category "components"
$ws.Cells.Item(70, 29).Value = "components"

# Update view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("AD70").Select()
